$d = $word.ActiveDocument

# Paragraph 1 currently reads:
#   "House of Hospitality, " <break> "Chapter Four ========================="
# Paragraph 2 currently reads (bold):
#   "By Dorothy Day"
#
# The edit removes paragraph 1 entirely and turns paragraph 2 into a plain
# (non-bold) pandoc-style title line: "% Dorothy Day"

$p1 = $d.Paragraphs(1)
$p1.Range.Delete()

$p2 = $d.Paragraphs(1)
$r = $p2.Range
$r.End = $r.End - 1        # exclude the paragraph mark so we don't touch pPr
$r.Delete()                 # drop the old "By Dorothy Day" text/formatting
$r.InsertAfter("% Dorothy Day")
